$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Range("C1").Value = "13-01-2023"

$data = @(
    @(2, "1810 Renta variable", 365309.07, 485699.02),
    @(3, "1822 Raices Valores Negociables", 819184.9, 819190.82),
    @(4, "Adcap IOL Acciones Argentina", 123274.27, 128213.71),
    @(5, "Allaria Acciones", 178433.93, 178472.19),
    @(6, "Alpha Acciones", 289342.64, 309685.18),
    @(7, "Alpha Latam", 123.92, 118.5),
    @(8, "Alpha Mega", 943595.25, 943783.88),
    @(9, "Alpha Mercosur", 625408.05, 660466.28),
    @(10, "Alpha Recursos Naturales", 311282.15, 342697.15),
    @(11, "Alpha planeam equil", 16532.89, 7610.35),
    @(12, "Alpha renta balan global", 2450230.98, 2449779.31),
    @(13, "Argenfunds", 39159.76, 39169.13),
    @(14, "Arpenta acciones", 8033.56, 8031.2),
    @(15, "Arpenta ex Mercosur", 14883.9, 14866.78),
    @(16, "Balanz", 628981.3, 638329.25),
    @(17, "Bull Market", 91964.45, 105811.67),
    @(18, "CMA acciones", 491519.94, 449341),
    @(19, "Compass Crecimiento", 2902240.16, 2882031.52),
    @(20, "Compass Crecimiento II", 17525.18, 17517.08),
    @(21, "Consultatio Acciones Argentina", 1537135.72, 1537159.6),
    @(22, "Consultatio Renta Variable", 607563.26, 607893.9),
    @(23, "Delta Acciones", 164274.08, 164346.87),
    @(24, "Delta Internacional", 6907.31, 6907.35),
    @(25, "Delta Latinoamerica", 6248.52, 6244.25),
    @(26, "Delta Recursos Naturales", 845439.2, 845910.37),
    @(27, "Delta Select", 1169260.18, 1158909.71),
    @(28, "Delta gestion V", 318428.81, 318069.43),
    @(29, "FBA Acciones Argentinas", 647904.57, 612802.31),
    @(30, "FBA Calificado", 638404.61, 603415.98),
    @(31, "Fima Acciones", 1055118.95, 1047369.39),
    @(32, "Fima PB Acciones", 694045.97, 632766.52),
    @(33, "Gainvest Renta Variable", 298166.4, 298315.05),
    @(34, "Galileo Acciones", 3776352.36, 3800311.45),
    @(35, "Goal Acciones Argentinas", 92432.95, 92453.92),
    @(36, "Goal acciones plus", 20703.83, 20696.69),
    @(37, "HF Acciones Argentinas", 413382.91, 403666.15),
    @(38, "HF Acciones Lideres", 739484.97, 744958.36),
    @(39, "IAM Renta Variable", 126595.33, 132610.88),
    @(40, "IEB Value", 18153.67, 18141.81),
    @(41, "Lombardi", 108967.3, 120003.84),
    @(42, "MAF", 126857.86, 126840.84),
    @(43, "Megainver", 113087.21, 113135.83),
    @(44, "Pellegrini Acciones", 253467.02, 293560.77),
    @(45, "Pionero Acciones", 371937.45, 371673.13),
    @(46, "Premier Renta Variable", 166670.66, 176528.12),
    @(47, "Quinquela Acciones", 362622.93, 362441.46),
    @(48, "Rofex 20 Renta Variable", 256133.99, 256093.93),
    @(49, "SBS Acciones Argentina", 1247721.79, 1248096.14),
    @(50, "Schroeder RV", 2797333.67, 2797583.94),
    @(51, "Supefondo RV", 3035007.23, 3138396.02),
    @(52, "Superfondo ", 791305.73, 792215.97),
    @(53, "Supergestion", 265311.28, 265070.5),
    @(54, "Toronto Trust Multimercado", 179215.94, 179222.74),
    @(55, "Toronto trust Argy", 103625.8, 103640.49),
    @(56, "avg", 623561.03, 627375.33),
    @(57, "total", 33672295.76, 33878267.73)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

